$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'66.181.49"
$ws.Range("E2").Value = "'  +1.63%  "

# Row 3
$ws.Range("D3").Value = "'3.555.02"
$ws.Range("E3").Value = "'  +5.00%  "

# Row 4
$ws.Range("E4").Value = "'  -0.26%  "

# Row 5
$ws.Range("D5").Value = "'606.36"
$ws.Range("E5").Value = "'  +2.55%  "

# Row 6
$ws.Range("D6").Value = "'145.15"
$ws.Range("E6").Value = "'  +3.15%  "

# Row 7
$ws.Range("D7").Value = "'3.554.03"
$ws.Range("E7").Value = "'  +4.81%  "

# Row 8
$ws.Range("E8").Value = "'  +0.24%  "

# Row 9
$ws.Range("D9").Value = "'0.486"
$ws.Range("E9").Value = "'  +4.15%  "

# Row 10
$ws.Range("D10").Value = "'0.136"
$ws.Range("E10").Value = "'  +2.00%  "

# Row 11
$ws.Range("D11").Value = "'8.01"
$ws.Range("E11").Value = "'  +2.26%  "

# Row 12
$ws.Range("E12").Value = "'  +2.09%  "

# Row 13
$ws.Range("D13").Value = "'4.121.22"
$ws.Range("E13").Value = "'  +3.51%  "

# Row 14
$ws.Range("D14").Value = "'0.0000208"
$ws.Range("E14").Value = "'  +5.25%  "

# Row 15
$ws.Range("D15").Value = "'30.14"
$ws.Range("E15").Value = "'  +2.08%  "

# Row 16
$ws.Range("D16").Value = "'3.560.85"
$ws.Range("E16").Value = "'  +4.71%  "

# Row 17
$ws.Range("D17").Value = "'66.307.51"
$ws.Range("E17").Value = "'  +1.99%  "

# Row 18
$ws.Range("E18").Value = "'  -0.66%  "

# Row 19
$ws.Range("D19").Value = "'11.46"
$ws.Range("E19").Value = "'  +11.27%  "

# Row 20
$ws.Range("D20").Value = "'6.21"
$ws.Range("E20").Value = "'  +2.10%  "

# Row 21
$ws.Range("D21").Value = "'14.98"
$ws.Range("E21").Value = "'  +2.00%  "

# Row 22
$ws.Range("D22").Value = "'430.48"
$ws.Range("E22").Value = "'  +3.73%  "

# Row 23
$ws.Range("D23").Value = "'0.609"
$ws.Range("E23").Value = "'  +5.32%  "

# Row 24
$ws.Range("D24").Value = "'78.69"
$ws.Range("E24").Value = "'  +2.00%  "

# Row 25
$ws.Range("D25").Value = "'3.700.92"
$ws.Range("E25").Value = "'  +4.79%  "

# Row 26
$ws.Range("E26").Value = "'  +0.00%  "

# Row 27
$ws.Range("D27").Value = "'0.0000119"
$ws.Range("E27").Value = "'  +8.87%  "

# Row 28
$ws.Range("D28").Value = "'2.52"
$ws.Range("E28").Value = "'  +4.81%  "

# Row 29
$ws.Range("D29").Value = "'8.04"
$ws.Range("E29").Value = "'  +3.47%  "

# Row 30
$ws.Range("D30").Value = "'9.22"
$ws.Range("E30").Value = "'  +0.70%  "

# Row 31
$ws.Range("E31").Value = "'  +0.03%  "

# Row 32
$ws.Range("B32").Value = "'Kaspa"
$ws.Range("C32").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").Value = "'0.161"
$ws.Range("E32").Value = "'  +0.75%  "

# Row 33
$ws.Range("B33").Value = "'Fetch.AI"
$ws.Range("C33").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.49"
$ws.Range("E33").Value = "'  +2.29%  "

# Row 34
$ws.Range("D34").Value = "'3.554.80"
$ws.Range("E34").Value = "'  +4.74%  "

# Row 35
$ws.Range("D35").Value = "'25.39"
$ws.Range("E35").Value = "'  +4.20%  "

# Row 36
$ws.Range("E36").Value = "'  -0.03%  "

# Row 37
$ws.Range("E37").Value = "'  +3.26%  "

# Row 38
$ws.Range("D38").Value = "'7.89"
$ws.Range("E38").Value = "'  +5.41%  "

# Row 39
$ws.Range("D39").Value = "'5.64"
$ws.Range("E39").Value = "'  +2.79%  "

# Row 40
$ws.Range("E40").Value = "'  -0.31%  "

# Row 41
$ws.Range("D41").Value = "'171.13"
$ws.Range("E41").Value = "'  +0.21%  "

# Row 42
$ws.Range("D42").Value = "'0.0856"
$ws.Range("E42").Value = "'  -0.48%  "

# Row 43
$ws.Range("D43").Value = "'5.20"
$ws.Range("E43").Value = "'  +3.61%  "

# Row 44
$ws.Range("D44").Value = "'0.897"
$ws.Range("E44").Value = "'  +3.81%  "

# Row 45
$ws.Range("D45").Value = "'1.95"
$ws.Range("E45").Value = "'  +2.32%  "

# Row 46
$ws.Range("D46").Value = "'46.07"
$ws.Range("E46").Value = "'  +1.72%  "

# Row 47
$ws.Range("E47").Value = "'  +4.26%  "

# Row 48
$ws.Range("D48").Value = "'25.88"
$ws.Range("E48").Value = "'  -1.93%  "

# Row 49
$ws.Range("E49").Value = "'  +4.88%  "

# Row 50
$ws.Range("D50").Value = "'7.13"
$ws.Range("E50").Value = "'  +1.17%  "

# Row 51
$ws.Range("D51").Value = "'0.955"
$ws.Range("E51").Value = "'  +4.80%  "
